# Update the "Förändrad" (Changed) date in column C for every data row
# (rows 2-255) from 2026-02-22 (46075) to 2026-02-23 (46076).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 255 }

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = 46076
}
